# Applies the edits described by the commit "Fixed a lot of errors"
$wb = $excel.ActiveWorkbook

# --- "main menu" sheet edits ---
$mainMenu = $wb.Worksheets.Item("main menu")

# Move "button must lead to a purchase and item detail screen" from D3 -> E3
$d3Value = $mainMenu.Range("D3").Value2
$mainMenu.Range("D3").Clear()
$mainMenu.Range("E3").Value2 = $d3Value

# Move "make button disabled ... nothing was searched" from C5 -> D5,
# fixing the typo "it" -> "if" at the same time
$mainMenu.Range("C5").Clear()
$mainMenu.Range("D5").Value2 = "make button disabled if nothing was searched"

# Update the active selection to E3
$mainMenu.Range("E3").Select()

# --- "account page" sheet edits ---
$accountPage = $wb.Worksheets.Item("account page")

# Move "display personnalized purchase history //" from C2 -> E2
$c2Value = $accountPage.Range("C2").Value2
$accountPage.Range("C2").Clear()
$accountPage.Range("E2").Value2 = $c2Value

# Fix the text in D2: "users" -> "only users"
$accountPage.Range("D2").Value2 = " display only users items currently for sale"

# "account page" was the originally active sheet; restore it as the active
# sheet/tab so that selecting a cell on "main menu" above does not change
# which sheet is active in the saved workbook.
$accountPage.Activate()
